$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Copy F1's eventual border/font style from the CURRENT E1 (which
#    already carries the "last column" header style: bold font + the
#    medium/thin outer-box border) before we touch E1 itself.
# ---------------------------------------------------------------------
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Re-point E1 to the shared "inner" header style (same as B1/C1/D1) now
# that it is no longer the right-most header cell.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# New header label for the appended category column.
$ws.Range("F1").Value = "id_categoria"

# ---------------------------------------------------------------------
# 2) Body rows: copy each row's existing look into the new F column so
#    the new cells inherit the matching thin/medium borders.
# ---------------------------------------------------------------------
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

$ws.Range("E3:E7").Copy()
$ws.Range("F3:F7").PasteSpecial(-4122)

# Fill in the new "id_categoria" values (all rows use category 2).
$ws.Range("F2:F7").Value = 2

# ---------------------------------------------------------------------
# 3) Data edits already present in the diff: id_funcionario goes from
#    2 to 3 for every data row.
# ---------------------------------------------------------------------
$ws.Range("A2:A7").Value = 3

# ---------------------------------------------------------------------
# 4) Leftover manual formatting noise that shipped in the same commit:
#    an underline highlight on the max price cell and an underline
#    "note" cell below the table, plus a stray normal-style touch.
# ---------------------------------------------------------------------
$ws.Range("C6").Font.Underline = 2
$ws.Range("B9").Font.Underline = 2
$ws.Range("F8").Font.Name = $ws.Range("F8").Font.Name

Write-Host "done"
